# Slide 4 ("Блокова схема" / block diagram): sensor + LCD labels updated.
$p = $ppt.ActivePresentation
$s4 = $p.Slides.Item(4)

# TextBox 28 (shape #12 on the slide): "AM2302" -> "DHT22"
$s4.Shapes.Item(12).TextFrame.TextRange.Text = "DHT22"

# TextBox 39 (shape #15 on the slide): "12x2 LCD" -> "16x2 LCD",
# splitting into two runs ("16x2 " + "LCD") the way PowerPoint does
# when only part of a run's text is retyped.
$tb39 = $s4.Shapes.Item(15).TextFrame.TextRange
$tb39.Characters(1, 5).Text = "16x2 "

# Slide 6 ("СПИСЪК С КОМПОНЕНТИ" / components list).
$s6 = $p.Slides.Item(6)
$content = $s6.Shapes.Item(2).TextFrame.TextRange

# Paragraph 2: "...AM2302" -> "...AM2302 (WIRED DHT22)"
$para2 = $content.Paragraphs(2)
$para2.Characters(33, 6).Text = "AM2302 (WIRED DHT22)"

# Paragraph 4: "16X2 " + "LCD" (two runs) -> "16X2 LCD" (single run)
$para4 = $content.Paragraphs(4)
$para4.Characters(3, 8).Text = "16X2 LCD"
